$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits: rename a few group labels in column B ---
# Header: "group_name" -> "cell_group"
$ws.Range("B1").Value = "cell_group"

# "philippines" -> "papua_new_guinea" (rows 20-22)
$ws.Range("B20:B22").Value = "papua_new_guinea"

# "orkney" -> "faeroe" (rows 41-44)
$ws.Range("B41:B44").Value = "faeroe"

# "andamans" -> "e_indian_ocean" (rows 48-51)
$ws.Range("B48:B51").Value = "e_indian_ocean"

# --- View / selection state ---
# Put the active selection/cell on B20:B22 (matches the saved selection in the sheet view)
$ws.Range("B20:B22").Select()

# Scroll the view so row 9 is at the top-left (matches topLeftCell="A9")
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1

# Reposition the application window (matches workbookView xWindow/yWindow)
$excel.Left = 10000
$excel.Top = 800
